$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 463.35715
$ws.Range("I33").Value = 307
$ws.Range("J33").Value = 619.7143
$ws.Range("K33").Value = 307
$ws.Range("L33").Value = 619.7143
$ws.Range("M33").Value = -78
$ws.Range("N33").Value = -1077.7143

$ws.Range("H101").Value = 493.5
$ws.Range("I101").Value = 470.75
$ws.Range("J101").Value = 584.5
$ws.Range("K101").Value = 1412.25
$ws.Range("L101").Value = 1753.5
$ws.Range("M101").Value = 209.75
$ws.Range("N101").Value = -4997.5

$ws.Range("H137").Value = 2854.2727
$ws.Range("J137").Value = 4099.75
$ws.Range("L137").Value = 12299.25
$ws.Range("N137").Value = -17399.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4883.2
$ws.Range("I3").Value = 3602.5
$ws.Range("K3").Value = 3602.5
$ws.Range("M3").Value = -3487.5

$ws.Range("H10").Value = 2900
$ws.Range("I10").Value = 2900
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2900
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -2730
$ws.Range("N10").Value = ""

$ws.Range("H12").Value = 5051.6665
$ws.Range("I12").Value = 201.5
$ws.Range("J12").Value = 7476.75
$ws.Range("K12").Value = 201.5
$ws.Range("L12").Value = 7476.75
$ws.Range("M12").Value = -28.5
$ws.Range("N12").Value = -7822.75

$ws.Range("H14").Value = 10459.6
$ws.Range("I14").Value = 17032.666
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 17032.666
$ws.Range("L14").Value = 600
$ws.Range("M14").Value = -16857.666
$ws.Range("N14").Value = -950

$ws.Range("H15").Value = 4611.3
$ws.Range("J15").Value = 4611.3
$ws.Range("L15").Value = 4611.3
$ws.Range("N15").Value = -5311.3

$ws.Range("H16").Value = 20607
$ws.Range("I16").Value = 33892.668
$ws.Range("J16").Value = 678.5
$ws.Range("K16").Value = 33892.668
$ws.Range("L16").Value = 678.5
$ws.Range("M16").Value = -33605.668
$ws.Range("N16").Value = -1252.5

$ws.Range("H17").Value = 879.5
$ws.Range("J17").Value = 879.5
$ws.Range("L17").Value = 879.5
$ws.Range("N17").Value = -1225.5

$ws.Range("H18").Value = 1256.5
$ws.Range("J18").Value = 1256.5
$ws.Range("L18").Value = 1256.5
$ws.Range("N18").Value = -1900.5

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = ""

$ws.Range("H21").Value = 1757.5
$ws.Range("I21").Value = 1757.5
$ws.Range("K21").Value = 1757.5
$ws.Range("M21").Value = -1383.5

$ws.Range("H32").Value = 9403
$ws.Range("I32").Value = 8115.6875
$ws.Range("K32").Value = 8115.6875
$ws.Range("M32").Value = -7828.6875

$ws.Range("H45").Value = 1867.2307
$ws.Range("I45").Value = 1856.1666
$ws.Range("K45").Value = 1856.1666
$ws.Range("M45").Value = -1479.1666

$ws.Range("H110").Value = 2343.6
$ws.Range("I110").Value = 583
$ws.Range("K110").Value = 583
$ws.Range("M110").Value = 1462

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 716.6667
$ws.Range("I80").Value = 250
$ws.Range("K80").Value = 250
$ws.Range("M80").Value = 748

$ws.Range("H83").Value = 716.6667
$ws.Range("I83").Value = 250
$ws.Range("K83").Value = 1250
$ws.Range("M83").Value = 3742

$ws.Range("H99").Value = 1331.3334
$ws.Range("I99").Value = 997.6
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 997.6
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = 500.4
$ws.Range("N99").Value = -5996

$ws.Range("H107").Value = 685.8
$ws.Range("I107").Value = 685.8
$ws.Range("K107").Value = 685.8
$ws.Range("M107").Value = 1234.2

$ws.Range("H134").Value = 7201.1724
$ws.Range("I134").Value = 7780.478
$ws.Range("K134").Value = 23341.434
$ws.Range("M134").Value = -20806.434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7853.846
$ws.Range("I58").Value = 4650
$ws.Range("K58").Value = 4650
$ws.Range("M58").Value = -4447

$ws.Range("H86").Value = 5137.4287
$ws.Range("I86").Value = 4792.4
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 4792.4
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -3669.4
$ws.Range("N86").Value = -8246

$ws.Range("H89").Value = 5137.4287
$ws.Range("I89").Value = 4792.4
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 23962
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -18346
$ws.Range("N89").Value = -41232

$ws.Range("H132").Value = 1304.5
$ws.Range("I132").Value = 1105
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 3315
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -785
$ws.Range("N132").Value = -15557

$ws.Range("H136").Value = 7853.846
$ws.Range("I136").Value = 4650
$ws.Range("K136").Value = 13950
$ws.Range("M136").Value = -11400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 88.23077000000001
$ws.Range("J12").Value = 98.7
$ws.Range("L12").Value = 296.1
$ws.Range("N12").Value = -642.1

$ws.Range("H38").Value = 37
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""

$ws.Range("H122").Value = 638.1429000000001
$ws.Range("I122").Value = 615.4
$ws.Range("K122").Value = 5538.599999999999
$ws.Range("M122").Value = -3088.599999999999

$ws.Range("H132").Value = 1429.1428
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2001.3334
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 18012.0006
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -23072.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 17862012
$ws.Range("I122").Value = 25004598
$ws.Range("K122").Value = 75013794
$ws.Range("M122").Value = -75011344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""

$ws.Range("H22").Value = 1682.5
$ws.Range("I22").Value = 1411.75
$ws.Range("J22").Value = 2224
$ws.Range("K22").Value = 1411.75
$ws.Range("L22").Value = 2224
$ws.Range("M22").Value = -1116.75
$ws.Range("N22").Value = -2814

$ws.Range("H27").Value = 1682.5
$ws.Range("I27").Value = 1411.75
$ws.Range("J27").Value = 2224
$ws.Range("K27").Value = 1411.75
$ws.Range("L27").Value = 2224
$ws.Range("M27").Value = -1304.75
$ws.Range("N27").Value = -2438

$ws.Range("H40").Value = 9165
$ws.Range("I40").Value = 7495
$ws.Range("K40").Value = 7495
$ws.Range("M40").Value = -7359

$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31498

$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -97488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""

$ws.Range("H107").Value = 863
$ws.Range("I107").Value = 863
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2589
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -669
$ws.Range("N107").Value = ""

$ws.Range("H136").Value = 2148.25
$ws.Range("I136").Value = 1991
$ws.Range("J136").Value = 3249
$ws.Range("K136").Value = 5973
$ws.Range("L136").Value = 9747
$ws.Range("M136").Value = -3423
$ws.Range("N136").Value = -14847
